$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting for new rows 17-19 (column A) to match existing style (bold, centered, bordered)
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null

# Column A sequential index values (row-2) for rows 10-19
for ($r = 10; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Row 10: Gaussian-Quadrature
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(10, 3).Value = 1.010334526185207
$ws.Cells.Item(10, 4).Value = 0.926759645655817
$ws.Cells.Item(10, 5).Value = 1.009172679661215
$ws.Cells.Item(10, 6).Value = 1.010334526185207
$ws.Cells.Item(10, 7).Value = 0.9599991728739697
$ws.Cells.Item(10, 8).Value = 1.03691429481915
$ws.Cells.Item(10, 9).Value = 1.012069983939075
$ws.Cells.Item(10, 10).Value = 0.926759645655817
$ws.Cells.Item(10, 11).Value = 0.9679661626585159
$ws.Cells.Item(10, 12).Value = 0.9891503444218617
$ws.Cells.Item(10, 13).Value = 0.9925417171890726

# Row 11: Spiral-90deg-10rot-5space
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11, 3).Value = 0.9918250382859393
$ws.Cells.Item(11, 4).Value = 1.007082951220033
$ws.Cells.Item(11, 5).Value = 0.9918111722493348
$ws.Cells.Item(11, 6).Value = 0.9918250382859393
$ws.Cells.Item(11, 7).Value = 1.000616530771101
$ws.Cells.Item(11, 8).Value = 0.9849033623358842
$ws.Cells.Item(11, 9).Value = 0.9912583617258226
$ws.Cells.Item(11, 10).Value = 1.007082951220033
$ws.Cells.Item(11, 11).Value = 0.999447061734684
$ws.Cells.Item(11, 12).Value = 0.9956360500103116
$ws.Cells.Item(11, 13).Value = 0.9945829027646859

# Row 12: Spiral-90deg-15rot-5space
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12, 3).Value = 0.991733114489273
$ws.Cells.Item(12, 4).Value = 1.007379881875671
$ws.Cells.Item(12, 5).Value = 0.9917567758855552
$ws.Cells.Item(12, 6).Value = 0.991733114489273
$ws.Cells.Item(12, 7).Value = 1.000751746359043
$ws.Cells.Item(12, 8).Value = 0.9847302702994051
$ws.Cells.Item(12, 9).Value = 0.9911916238573156
$ws.Cells.Item(12, 10).Value = 1.007379881875671
$ws.Cells.Item(12, 11).Value = 0.9995683288806133
$ws.Cells.Item(12, 12).Value = 0.9956507216849433
$ws.Cells.Item(12, 13).Value = 0.9945905687943771

# Row 13: Spiral-90deg-10rot-3space
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13, 3).Value = 0.9918159117793458
$ws.Cells.Item(13, 4).Value = 1.007104857164515
$ws.Cells.Item(13, 5).Value = 0.9918197796544369
$ws.Cells.Item(13, 6).Value = 0.9918159117793458
$ws.Cells.Item(13, 7).Value = 1.00061857131567
$ws.Cells.Item(13, 8).Value = 0.9848648329629496
$ws.Cells.Item(13, 9).Value = 0.9912539182258197
$ws.Cells.Item(13, 10).Value = 1.007104857164515
$ws.Cells.Item(13, 11).Value = 0.999462318409476
$ws.Cells.Item(13, 12).Value = 0.9956391150944109
$ws.Cells.Item(13, 13).Value = 0.9945796451837895

# Row 14: NoRotation-tilt60deg
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(14, 3).Value = 1.020635999999999
$ws.Cells.Item(14, 4).Value = 0.8703680000000006
$ws.Cells.Item(14, 5).Value = 1.027928
$ws.Cells.Item(14, 6).Value = 1.020635999999999
$ws.Cells.Item(14, 7).Value = 0.922472000000002
$ws.Cells.Item(14, 8).Value = 1.097543999999999
$ws.Cells.Item(14, 9).Value = 1.026679999999999
$ws.Cells.Item(14, 10).Value = 0.8703680000000006
$ws.Cells.Item(14, 11).Value = 0.9491480000000003
$ws.Cells.Item(14, 12).Value = 0.9848919999999997
$ws.Cells.Item(14, 13).Value = 0.9942713333333332

# Row 15: Rotation-NoTilt
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Cells.Item(15, 3).Value = 1.05
$ws.Cells.Item(15, 4).Value = 0.7461125000000011
$ws.Cells.Item(15, 5).Value = 1.06
$ws.Cells.Item(15, 6).Value = 1.05
$ws.Cells.Item(15, 7).Value = 0.85
$ws.Cells.Item(15, 8).Value = 1.2
$ws.Cells.Item(15, 9).Value = 1.06
$ws.Cells.Item(15, 10).Value = 0.7461125000000011
$ws.Cells.Item(15, 11).Value = 0.9030562500000006
$ws.Cells.Item(15, 12).Value = 0.9765281250000002
$ws.Cells.Item(15, 13).Value = 0.9943520833333336

# Row 16: Rotation-60detTilt
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"
$ws.Cells.Item(16, 3).Value = 1.027505905049599
$ws.Cells.Item(16, 4).Value = 0.8493880448000025
$ws.Cells.Item(16, 5).Value = 1.033084348825599
$ws.Cells.Item(16, 6).Value = 1.027505905049599
$ws.Cells.Item(16, 7).Value = 0.9105357925376051
$ws.Cells.Item(16, 8).Value = 1.112871298457598
$ws.Cells.Item(16, 9).Value = 1.032606124441603
$ws.Cells.Item(16, 10).Value = 0.8493880448000025
$ws.Cells.Item(16, 11).Value = 0.9412361968128009
$ws.Cells.Item(16, 12).Value = 0.9843710509311998
$ws.Cells.Item(16, 13).Value = 0.9943319190186678

# Row 17: HexGrid-90degTilt5degRes
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17, 3).Value = 0.9946196628895481
$ws.Cells.Item(17, 4).Value = 0.99512277449894
$ws.Cells.Item(17, 5).Value = 0.9953239227158313
$ws.Cells.Item(17, 6).Value = 0.9946196628895481
$ws.Cells.Item(17, 7).Value = 0.9944352433466371
$ws.Cells.Item(17, 8).Value = 0.995624102483925
$ws.Cells.Item(17, 9).Value = 0.9949609971634735
$ws.Cells.Item(17, 10).Value = 0.99512277449894
$ws.Cells.Item(17, 11).Value = 0.9952233486073856
$ws.Cells.Item(17, 12).Value = 0.9949215057484668
$ws.Cells.Item(17, 13).Value = 0.9950144505163925

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18, 3).Value = 0.9903470099055087
$ws.Cells.Item(18, 4).Value = 1.003349142013436
$ws.Cells.Item(18, 5).Value = 0.9928917274737754
$ws.Cells.Item(18, 6).Value = 0.9903470099055087
$ws.Cells.Item(18, 7).Value = 0.9995116746476891
$ws.Cells.Item(18, 8).Value = 0.9918743210785456
$ws.Cells.Item(18, 9).Value = 0.993469216594414
$ws.Cells.Item(18, 10).Value = 1.003349142013436
$ws.Cells.Item(18, 11).Value = 0.9981204347436059
$ws.Cells.Item(18, 12).Value = 0.9942337223245573
$ws.Cells.Item(18, 13).Value = 0.9952405152855613

# Row 19: HexGrid-60degTilt5degRes
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19, 3).Value = 0.9885419359001381
$ws.Cells.Item(19, 4).Value = 1.023904868801994
$ws.Cells.Item(19, 5).Value = 0.988328431585138
$ws.Cells.Item(19, 6).Value = 0.9885419359001381
$ws.Cells.Item(19, 7).Value = 1.011093490859023
$ws.Cells.Item(19, 8).Value = 0.9748694217102131
$ws.Cells.Item(19, 9).Value = 0.9874894006508635
$ws.Cells.Item(19, 10).Value = 1.023904868801994
$ws.Cells.Item(19, 11).Value = 1.006116650193566
$ws.Cells.Item(19, 12).Value = 0.9973292930468519
$ws.Cells.Item(19, 13).Value = 0.9957045915845617

Write-Output "done"